$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing Science score for "hari" (row 6), which was previously blank.
$ws.Range("C6").Value = 85

# Recompute Total / Average formulas for all student rows so the newly
# entered value is reflected (also re-applies the formulas individually,
# matching the per-row Total/Average formulas in the sheet).
$ws.Range("G3").Formula = "=SUM(C3:F3)"
$ws.Range("H3").Formula = "=AVERAGE(C3:F3)"
$ws.Range("G4").Formula = "=SUM(C4:F4)"
$ws.Range("H4").Formula = "=AVERAGE(C4:F4)"
$ws.Range("G5").Formula = "=SUM(C5:F5)"
$ws.Range("H5").Formula = "=AVERAGE(C5:F5)"
$ws.Range("G6").Formula = "=SUM(C6:F6)"
$ws.Range("H6").Formula = "=AVERAGE(C6:F6)"

# Move the active selection to J7, matching where the user clicked next.
[void]$ws.Range("J7").Select()
